# Update branch-wise NIL/US/SS summary figures on the active worksheet.
# This refreshes the underlying counts that drive the report's layout/table
# (design layout + table style refresh), row by row, column by column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 34
$ws.Cells.Item(2, 4).Value = 10
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 16
$ws.Cells.Item(2, 7).Value = 38
$ws.Cells.Item(3, 3).Value = 38
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 7).Value = 31
$ws.Cells.Item(4, 3).Value = 37
$ws.Cells.Item(4, 5).Value = 12
$ws.Cells.Item(4, 7).Value = 27
$ws.Cells.Item(5, 3).Value = 35
$ws.Cells.Item(5, 4).Value = 14
$ws.Cells.Item(5, 6).Value = 13
$ws.Cells.Item(5, 7).Value = 41
$ws.Cells.Item(6, 3).Value = 32
$ws.Cells.Item(6, 4).Value = 13
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 18
$ws.Cells.Item(6, 7).Value = 38
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 31
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = 8
$ws.Cells.Item(7, 7).Value = 45
$ws.Cells.Item(8, 3).Value = 32
$ws.Cells.Item(8, 4).Value = 15
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 14
$ws.Cells.Item(8, 7).Value = 40
$ws.Cells.Item(9, 3).Value = 30
$ws.Cells.Item(9, 4).Value = 8
$ws.Cells.Item(9, 5).Value = 6
$ws.Cells.Item(9, 6).Value = 10
$ws.Cells.Item(9, 7).Value = 52
$ws.Cells.Item(10, 3).Value = 39
$ws.Cells.Item(10, 4).Value = 13
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = 10
$ws.Cells.Item(10, 7).Value = 39
$ws.Cells.Item(11, 3).Value = 32
$ws.Cells.Item(11, 4).Value = 14
$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 8
$ws.Cells.Item(11, 7).Value = 42
$ws.Cells.Item(12, 3).Value = 34
$ws.Cells.Item(12, 4).Value = 20
$ws.Cells.Item(12, 5).Value = 9
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(12, 7).Value = 38
$ws.Cells.Item(13, 3).Value = 36
$ws.Cells.Item(13, 4).Value = 27
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 10
$ws.Cells.Item(13, 7).Value = 24
$ws.Cells.Item(14, 3).Value = 35
$ws.Cells.Item(14, 4).Value = 16
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 7).Value = 41
$ws.Cells.Item(15, 3).Value = 30
$ws.Cells.Item(15, 4).Value = 6
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 8
$ws.Cells.Item(15, 7).Value = 52
$ws.Cells.Item(16, 3).Value = 35
$ws.Cells.Item(16, 4).Value = 11
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 15
$ws.Cells.Item(16, 7).Value = 36
$ws.Cells.Item(17, 4).Value = 25
$ws.Cells.Item(17, 6).Value = 4
$ws.Cells.Item(17, 7).Value = 26
$ws.Cells.Item(18, 3).Value = 33
$ws.Cells.Item(18, 4).Value = 7
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 7).Value = 59
$ws.Cells.Item(19, 3).Value = 32
$ws.Cells.Item(19, 4).Value = 10
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 5
$ws.Cells.Item(19, 7).Value = 56
$ws.Cells.Item(20, 3).Value = 38
$ws.Cells.Item(20, 4).Value = 15
$ws.Cells.Item(20, 5).Value = 12
$ws.Cells.Item(20, 6).Value = 11
$ws.Cells.Item(20, 7).Value = 30
$ws.Cells.Item(21, 3).Value = 37
$ws.Cells.Item(21, 4).Value = 11
$ws.Cells.Item(21, 5).Value = 7
$ws.Cells.Item(22, 3).Value = 36
$ws.Cells.Item(22, 4).Value = 24
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 12
$ws.Cells.Item(22, 7).Value = 31
$ws.Cells.Item(23, 4).Value = 11
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = 12
$ws.Cells.Item(23, 7).Value = 28
$ws.Cells.Item(24, 3).Value = 36
$ws.Cells.Item(24, 4).Value = 23
$ws.Cells.Item(24, 5).Value = 5
$ws.Cells.Item(24, 6).Value = 5
$ws.Cells.Item(24, 7).Value = 37
$ws.Cells.Item(25, 3).Value = 32
$ws.Cells.Item(25, 5).Value = 5
$ws.Cells.Item(25, 6).Value = 14
$ws.Cells.Item(25, 7).Value = 50
$ws.Cells.Item(26, 3).Value = 31
$ws.Cells.Item(26, 4).Value = 10
$ws.Cells.Item(26, 5).Value = 13
$ws.Cells.Item(26, 7).Value = 40
$ws.Cells.Item(27, 4).Value = 14
$ws.Cells.Item(27, 5).Value = 7
$ws.Cells.Item(27, 6).Value = 9
$ws.Cells.Item(27, 7).Value = 46
$ws.Cells.Item(28, 3).Value = 38
$ws.Cells.Item(28, 4).Value = 11
$ws.Cells.Item(28, 5).Value = 6
$ws.Cells.Item(28, 6).Value = 11
$ws.Cells.Item(28, 7).Value = 40
$ws.Cells.Item(29, 3).Value = 33
$ws.Cells.Item(29, 5).Value = 6
$ws.Cells.Item(30, 3).Value = 34
$ws.Cells.Item(30, 4).Value = 17
$ws.Cells.Item(30, 5).Value = 11
$ws.Cells.Item(30, 6).Value = 14
$ws.Cells.Item(30, 7).Value = 30
$ws.Cells.Item(31, 3).Value = 34
$ws.Cells.Item(31, 4).Value = 22
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(31, 6).Value = 12
$ws.Cells.Item(31, 7).Value = 30
$ws.Cells.Item(32, 3).Value = 35
$ws.Cells.Item(32, 4).Value = 17
$ws.Cells.Item(32, 5).Value = 6
$ws.Cells.Item(32, 7).Value = 42
